# Applies updated Leve profit calculations (currentAveragePrice / LevePrice / LeveProfit
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets, matching the
# scheduled-runner refresh of the Kujata Profits workbook.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3638.2173
$ws.Range("I76").Value = 3810.5881
$ws.Range("K76").Value = 3810.5881
$ws.Range("M76").Value = -3495.5881
$ws.Range("H79").Value = 3638.2173
$ws.Range("I79").Value = 3810.5881
$ws.Range("K79").Value = 3810.5881
$ws.Range("M79").Value = -2718.5881
$ws.Range("H117").Value = 29999
$ws.Range("J117").Value = 29999
$ws.Range("L117").Value = 29999
$ws.Range("N117").Value = -39177
$ws.Range("H139").Value = 31948.889
$ws.Range("J139").Value = 31948.889
$ws.Range("L139").Value = 31948.889
$ws.Range("N139").Value = -42228.889
$ws.Range("H141").Value = 848.3333
$ws.Range("I141").Value = 743.63635
$ws.Range("K141").Value = 2230.90905
$ws.Range("M141").Value = 2949.09095

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 498.33334
$ws.Range("I26").Value = 303.5
$ws.Range("J26").Value = 888
$ws.Range("K26").Value = 303.5
$ws.Range("L26").Value = 888
$ws.Range("M26").Value = 26.5
$ws.Range("N26").Value = -1548
$ws.Range("H32").Value = 4751.1665
$ws.Range("I32").Value = 4344.894
$ws.Range("K32").Value = 4344.894
$ws.Range("M32").Value = -4057.894
$ws.Range("H36").Value = 10754.8
$ws.Range("I36").Value = 10754.8
$ws.Range("K36").Value = 10754.8
$ws.Range("M36").Value = -10408.8
$ws.Range("H38").Value = 7000
$ws.Range("I38").Value = 7000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 7000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -6533
$ws.Range("N38").Value = $null
$ws.Range("H74").Value = 1183.5143
$ws.Range("I74").Value = 920.3333
$ws.Range("J74").Value = 2762.6
$ws.Range("K74").Value = 920.3333
$ws.Range("L74").Value = 2762.6
$ws.Range("M74").Value = -46.33330000000001
$ws.Range("N74").Value = -4510.6
$ws.Range("H77").Value = 1183.5143
$ws.Range("I77").Value = 920.3333
$ws.Range("J77").Value = 2762.6
$ws.Range("K77").Value = 4601.6665
$ws.Range("L77").Value = 13813
$ws.Range("M77").Value = -233.6665000000003
$ws.Range("N77").Value = -22549
$ws.Range("H92").Value = 2750000
$ws.Range("J92").Value = 2750000
$ws.Range("L92").Value = 2750000
$ws.Range("N92").Value = -2754992

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 179.33333
$ws.Range("I7").Value = 179.33333
$ws.Range("K7").Value = 179.33333
$ws.Range("M7").Value = -66.33332999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 5000308
$ws.Range("I12").Value = 400.5
$ws.Range("J12").Value = 10000215
$ws.Range("K12").Value = 400.5
$ws.Range("L12").Value = 10000215
$ws.Range("M12").Value = -230.5
$ws.Range("N12").Value = -10000555
$ws.Range("H31").Value = 2065.72
$ws.Range("J31").Value = 1356
$ws.Range("L31").Value = 1356
$ws.Range("N31").Value = -1946
$ws.Range("H34").Value = 2065.72
$ws.Range("J34").Value = 1356
$ws.Range("L34").Value = 1356
$ws.Range("N34").Value = -1760
$ws.Range("H35").Value = 850
$ws.Range("I35").Value = 850
$ws.Range("K35").Value = 850
$ws.Range("M35").Value = -556

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 122.71429
$ws.Range("I2").Value = 69.666664
$ws.Range("J2").Value = 162.5
$ws.Range("K2").Value = 417.999984
$ws.Range("L2").Value = 975
$ws.Range("M2").Value = -304.999984
$ws.Range("N2").Value = -1201
$ws.Range("H4").Value = 450488.78
$ws.Range("I4").Value = 50018.277
$ws.Range("J4").Value = 1351547.4
$ws.Range("K4").Value = 150054.831
$ws.Range("L4").Value = 4054642.2
$ws.Range("M4").Value = -149942.831
$ws.Range("N4").Value = -4054866.2
$ws.Range("H134").Value = 3144.0303
$ws.Range("I134").Value = 1189.8948
$ws.Range("J134").Value = 5796.0713
$ws.Range("K134").Value = 3569.6844
$ws.Range("L134").Value = 17388.2139
$ws.Range("M134").Value = 1500.3156
$ws.Range("N134").Value = -27528.2139

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1531.45
$ws.Range("J7").Value = 1501.6666
$ws.Range("L7").Value = 1501.6666
$ws.Range("N7").Value = -1725.6666
$ws.Range("H22").Value = 694.2222
$ws.Range("I22").Value = 611.5
$ws.Range("K22").Value = 611.5
$ws.Range("M22").Value = -316.5
$ws.Range("H27").Value = 694.2222
$ws.Range("I27").Value = 611.5
$ws.Range("K27").Value = 611.5
$ws.Range("M27").Value = -504.5
$ws.Range("H64").Value = 21575
$ws.Range("J64").Value = 21575
$ws.Range("L64").Value = 21575
$ws.Range("N64").Value = -22025
$ws.Range("H67").Value = 21575
$ws.Range("J67").Value = 21575
$ws.Range("L67").Value = 21575
$ws.Range("N67").Value = -23135
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352
$ws.Range("H126").Value = 1531.45
$ws.Range("J126").Value = 1501.6666
$ws.Range("L126").Value = 4504.9998
$ws.Range("N126").Value = -9444.9998
$ws.Range("H127").Value = 37500
$ws.Range("J127").Value = 37500
$ws.Range("L127").Value = 37500
$ws.Range("N127").Value = -47420
$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960
$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = $null
$ws.Range("H14").Value = 211254.75
$ws.Range("I14").Value = 352504.5
$ws.Range("K14").Value = 352504.5
$ws.Range("M14").Value = -352336.5
$ws.Range("H63").Value = 11218.667
$ws.Range("I63").Value = 4113
$ws.Range("J63").Value = 13248.857
$ws.Range("K63").Value = 4113
$ws.Range("L63").Value = 13248.857
$ws.Range("M63").Value = -3489
$ws.Range("N63").Value = -14496.857
$ws.Range("H66").Value = 11218.667
$ws.Range("I66").Value = 4113
$ws.Range("J66").Value = 13248.857
$ws.Range("K66").Value = 12339
$ws.Range("L66").Value = 39746.571
$ws.Range("M66").Value = -9219
$ws.Range("N66").Value = -45986.571
$ws.Range("H81").Value = 401
$ws.Range("I81").Value = 400.5
$ws.Range("J81").Value = 402
$ws.Range("K81").Value = 801
$ws.Range("L81").Value = 804
$ws.Range("M81").Value = 260
$ws.Range("N81").Value = -2926
$ws.Range("H82").Value = 14000
$ws.Range("J82").Value = 14000
$ws.Range("L82").Value = 14000
$ws.Range("N82").Value = -14766
$ws.Range("H84").Value = 401
$ws.Range("I84").Value = 400.5
$ws.Range("J84").Value = 402
$ws.Range("K84").Value = 4005
$ws.Range("L84").Value = 4020
$ws.Range("M84").Value = 1299
$ws.Range("N84").Value = -14628
$ws.Range("H85").Value = 14000
$ws.Range("J85").Value = 14000
$ws.Range("L85").Value = 14000
$ws.Range("N85").Value = -14000
